$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2026-02-08 Sunday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-02-09 Monday", 2)

# Update the multiplication table (single table, 5 columns; data in rows 1,5,10,15,20)
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text  = "15×92="
$t.Cell(1, 2).Range.Text  = "59×78="
$t.Cell(1, 3).Range.Text  = "99×39="
$t.Cell(1, 4).Range.Text  = "86×45="
$t.Cell(1, 5).Range.Text  = "67×84="

$t.Cell(5, 1).Range.Text  = "56×71="
$t.Cell(5, 2).Range.Text  = "47×69="
$t.Cell(5, 3).Range.Text  = "55×43="
$t.Cell(5, 4).Range.Text  = "90×27="
$t.Cell(5, 5).Range.Text  = "39×96="

$t.Cell(10, 1).Range.Text = "82×93="
$t.Cell(10, 2).Range.Text = "68×25="
$t.Cell(10, 3).Range.Text = "71×28="
$t.Cell(10, 4).Range.Text = "32×16="
$t.Cell(10, 5).Range.Text = "98×33="

$t.Cell(15, 1).Range.Text = "70×23="
$t.Cell(15, 2).Range.Text = "85×96="
$t.Cell(15, 3).Range.Text = "79×25="
$t.Cell(15, 4).Range.Text = "64×78="
$t.Cell(15, 5).Range.Text = "21×33="

$t.Cell(20, 1).Range.Text = "70×77="
$t.Cell(20, 2).Range.Text = "98×82="
$t.Cell(20, 3).Range.Text = "53×63="
$t.Cell(20, 4).Range.Text = "48×83="
$t.Cell(20, 5).Range.Text = "53×49="
